$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: fix calculated/"trash" values that had been stored as text strings
$ws.Range("M5").Value = 11538.46153846154
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("Q5").Value = 465.67

# Row 6: table fix - recalculated figures
$ws.Range("M6").Value = 11538.46153846154
$ws.Range("Q6").Value = 620.89
